$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Title) to make room for the SJR score,
# shifting Title / H index / Publisher / Categories one column to the right.
$ws.Columns("B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "SJR"

# SJR (SCImago Journal Rank) score for each of the 15 data rows.
$sjr = @(1.232, 1.177, 1.17, 0.954, 0.688, 0.532, 0.521, 0.464, 0.44, 0.298, 0.291, 0.261, 0.225, 0.208, 0.109)

for ($i = 0; $i -lt $sjr.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $sjr[$i]
}

Write-Output "done"
